$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "nu of poles" column T (one new data point per row, 2..32)
$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 2
    11 = 2
    12 = 2
    13 = 2
    14 = 2
    15 = 2
    16 = 2
    17 = 3
    18 = 3
    19 = 2
    20 = 3
    21 = 4
    22 = 3
    23 = 4
    24 = 4
    25 = 4
    26 = 4
    27 = 4
    28 = 4
    29 = 4
    30 = 4
    31 = 4
    32 = 5
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 20).Value = $values[$row]
}

# Match the author's post-edit view: scrolled right and the freshly
# populated next column (U, the now-empty column right after the new T
# data) selected top-to-bottom.
$excel.ActiveWindow.ScrollColumn = 18
$ws.Range("U2:U32").Select()
